$d = $word.ActiveDocument

# --- First paragraph: token text, indent, and border updates ---
$p1 = $d.Paragraphs(1)

# Replace the placeholder token text in the first run.
$d.Content.Find.Execute("**ID__AFFARS_5301_topic_27__ID**", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5301_7__ID**", 2)

# Remove the now-orphaned trailing space (formerly its own run) that
# followed the token, leaving a single run with just the token text.
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$trailing = $d.Range($r1.End - 2, $r1.End - 1)
if ($trailing.Text -eq " ") {
    $trailing.Delete()
}

# Update the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Add a paragraph border (top/left/bottom/right) with 5pt padding and no
# line, matching the target pBdr of <w:top w:space="5"/> etc.
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
